$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6 (this shifts the existing rows 6-31 down to 7-32,
# and copies the formatting of row 5 into the new row 6).
$ws.Rows("6:6").Insert()

# Column A of the new row should use the same format as column C of the Users
# block above (fillId2, general alignment) rather than the left-aligned format
# that Insert() copied down from A5. Copy just the formatting over.
$ws.Range("C2").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# Fill in the new "deleteUser" endpoint row.
$ws.Range("A6").Value = "deleteUser"
$ws.Range("B6").Value = "DELETE"
$ws.Range("C6").Value = "/user/:uid"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "200, 400, 404"
$ws.Range("F6").Value = "Object (User)"

# The sheet previously ended with three blank spacer rows (styled, no data).
# After inserting the row above, the first of those spacer rows (now row 30)
# is redundant, so remove it to keep the same overall row count / dimension.
$ws.Rows("30:30").Delete()

# Leave the sheet with the whole grid selected (matches the saved selection
# state of the edited workbook).
$ws.Cells.Select()
